$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Columns D (Price) and E (Volume(1h))
# contain numeric-looking text that Excel would otherwise auto-convert to
# real numbers/percentages, so we force text entry with a leading apostrophe
# and then clear the resulting quote-prefix formatting so the cell keeps the
# workbook's original (unstyled) look while remaining plain text.
$updates = @(
    @{ Cell = "D2"; Value = "'325.98" }
    @{ Cell = "E2"; Value = "'-2.88%" }
    @{ Cell = "D3"; Value = "'44.58" }
    @{ Cell = "E3"; Value = "'1.49%" }
    @{ Cell = "D4"; Value = "'5.560" }
    @{ Cell = "E4"; Value = "'-4.23%" }
    @{ Cell = "D5"; Value = "'0.08065" }
    @{ Cell = "E5"; Value = "'-3.31%" }
    @{ Cell = "D6"; Value = "'8.677" }
    @{ Cell = "E6"; Value = "'-1.73%" }
    @{ Cell = "D7"; Value = "'1.907" }
    @{ Cell = "E7"; Value = "'-3.96%" }
    @{ Cell = "D8"; Value = "'4.298" }
    @{ Cell = "E8"; Value = "'-4.70%" }
    @{ Cell = "D9"; Value = "'2.700" }
    @{ Cell = "E9"; Value = "'-6.83%" }
    @{ Cell = "D10"; Value = "'0.9418" }
    @{ Cell = "E10"; Value = "'0.12%" }
    @{ Cell = "E11"; Value = "'-5.44%" }
    @{ Cell = "E12"; Value = "'-4.77%" }
    @{ Cell = "D13"; Value = "'0.09962" }
    @{ Cell = "E13"; Value = "'0.82%" }
    @{ Cell = "D14"; Value = "'0.04267" }
    @{ Cell = "E14"; Value = "'-6.86%" }
    @{ Cell = "D15"; Value = "'0.1065" }
    @{ Cell = "E15"; Value = "'-0.29%" }
    @{ Cell = "D16"; Value = "'0.001279" }
    @{ Cell = "E16"; Value = "'-1.61%" }
    @{ Cell = "D17"; Value = "'0.04195" }
    @{ Cell = "E17"; Value = "'-4.64%" }
    @{ Cell = "D18"; Value = "'0.005855" }
    @{ Cell = "E18"; Value = "'-1.64%" }
    @{ Cell = "B19"; Value = "LEO" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D19"; Value = "'3.580" }
    @{ Cell = "E19"; Value = "'2.38%" }
    @{ Cell = "B20"; Value = "BitpandaEcosystemToken" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Cell = "D20"; Value = "'0.3503" }
    @{ Cell = "E20"; Value = "'-0.19%" }
    @{ Cell = "B21"; Value = "MCDex" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D21"; Value = "'8.397" }
    @{ Cell = "E21"; Value = "'-4.39%" }
    @{ Cell = "B22"; Value = "ProBitToken" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" }
    @{ Cell = "D22"; Value = "'0.1371" }
    @{ Cell = "E22"; Value = "'0.58%" }
    @{ Cell = "B23"; Value = "ZBToken" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb" }
    @{ Cell = "D23"; Value = "'0.2526" }
    @{ Cell = "E23"; Value = "'-3.28%" }
    @{ Cell = "B24"; Value = "BitKan" }
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan" }
    @{ Cell = "D24"; Value = "'0.001242" }
    @{ Cell = "E24"; Value = "'-1.32%" }
    @{ Cell = "B25"; Value = "HotbitToken" }
    @{ Cell = "C25"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" }
    @{ Cell = "D25"; Value = "'0.004495" }
    @{ Cell = "E25"; Value = "'2.29%" }
    @{ Cell = "E26"; Value = "'-6.44%" }
    @{ Cell = "D27"; Value = "'0.0003992" }
    @{ Cell = "E27"; Value = "'-0.04%" }
    @{ Cell = "D39"; Value = "'0.02637" }
    @{ Cell = "E39"; Value = "'-5.68%" }
    @{ Cell = "E40"; Value = "'-4.56%" }
    @{ Cell = "D41"; Value = "'0.007652" }
    @{ Cell = "E41"; Value = "'-3.40%" }
    @{ Cell = "D42"; Value = "'0.1398" }
    @{ Cell = "E42"; Value = "'-2.25%" }
    @{ Cell = "D43"; Value = "'0.006955" }
    @{ Cell = "E43"; Value = "'-22.43%" }
    @{ Cell = "D44"; Value = "'0.002024" }
    @{ Cell = "E44"; Value = "'-3.72%" }
    @{ Cell = "D45"; Value = "'0.008858" }
    @{ Cell = "E45"; Value = "'-16.32%" }
    @{ Cell = "D46"; Value = "'0.00007177" }
    @{ Cell = "E46"; Value = "'-1.55%" }
    @{ Cell = "D47"; Value = "'0.00000000751" }
    @{ Cell = "E47"; Value = "'0.00%" }
    @{ Cell = "D48"; Value = "'0.003535" }
    @{ Cell = "E48"; Value = "'9.00%" }
    @{ Cell = "D49"; Value = "'0.002272" }
    @{ Cell = "E49"; Value = "'-0.04%" }
    @{ Cell = "D50"; Value = "'0.00002102" }
    @{ Cell = "E50"; Value = "'0.00%" }
    @{ Cell = "D51"; Value = "'0.0002002" }
    @{ Cell = "E51"; Value = "'0.00%" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = $u.Value
    $range.ClearFormats()
}
